# Auto-generated edit script: updates crypto price/volume data per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.093.14"
$ws.Range("E2").Value = "  +4.79%  "
$ws.Range("D3").Value = "2.437.00"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.516"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +10.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.123"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.01"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("D15").Value = "2.818.58"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "2.436.54"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.76%  "
$ws.Range("D18").Value = "45.004.32"
$ws.Range("E18").Value = "  +4.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  +16.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +13.37%  "
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0765"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("E38").Value = "  +4.88%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.24"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("D45").Value = "1.948.75"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  +8.39%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.28"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +19.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.09%  "
